$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.419.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.200.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.77%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.197.11"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.543"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.72%  "
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.732.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.098.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.202.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.773"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.19%  "
$ws.Range("E23").Value = "  +6.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.09%  "
$ws.Range("E25").Value = "  +11.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.37%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +9.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0901"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.85%  "
$ws.Range("E37").Value = "  +4.75%  "
$ws.Range("E38").Value = "  +5.83%  "
$ws.Range("E39").Value = "  +3.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "480.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "52.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.97%  "
$ws.Range("E43").Value = "  +10.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0383"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.951.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("E46").Value = "  +4.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.34%  "
$ws.Range("E51").Value = "  -0.01%  "
